$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.249.88"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.907.51"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "307.84"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "0.3815"
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").Value = "0.07309"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("D10").Value = "21.63"
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("D11").Value = "0.9058"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "0.08171"
$ws.Range("E12").Value = "  -3.98%  "
$ws.Range("D13").Value = "96.47"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").Value = "5.373"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "1.666.71"
$ws.Range("E15").Value = "  -12.48%  "
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "0.000008691"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "14.73"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "27.282.28"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "5.130"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").Value = "10.82"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").Value = "6.512"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").Value = "2.340"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("D25").Value = "149.90"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "1.742"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").Value = "116.77"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").Value = "4.851"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("D30").Value = "4.879"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").Value = "0.09256"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "0.8246"
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("D33").Value = "0.05085"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("D34").Value = "1.228"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("D35").Value = "3.002"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "2.759"
$ws.Range("E36").Value = "  +5.16%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "3.363"
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").Value = "0.5763"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").Value = "0.02006"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").Value = "1.083"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").Value = "9.101"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").Value = "6.604"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").Value = "117.08"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "0.1524"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "0.4926"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.17"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "1.644"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").Value = "38.72"
$ws.Range("E49").Value = "  +3.06%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "64.02"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06052"
$ws.Range("E51").Value = "  +1.88%  "
